$wb = $excel.ActiveWorkbook
$after = $wb.Worksheets.Item("B2")
$ws = $wb.Worksheets.Add($null, $after)
$ws.Name = "B3"
$ws.Activate()
Write-Host ($excel.ActiveWindow | Get-Member -MemberType Property | Out-String)
